$wb = $excel.ActiveWorkbook

# ALC row 53
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 1123.1428
$ws.Range("J53").Value = 1374.4
$ws.Range("L53").Value = 1374.4
$ws.Range("N53").Value = -2648.4

# ALC row 70
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1471.8889
$ws.Range("I70").Value = 1099.5
$ws.Range("J70").Value = 1578.2858
$ws.Range("K70").Value = 3298.5
$ws.Range("L70").Value = 4734.857400000001
$ws.Range("M70").Value = -3028.5
$ws.Range("N70").Value = -5274.857400000001

# ALC row 73
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 1471.8889
$ws.Range("I73").Value = 1099.5
$ws.Range("J73").Value = 1578.2858
$ws.Range("K73").Value = 3298.5
$ws.Range("L73").Value = 4734.857400000001
$ws.Range("M73").Value = -2362.5
$ws.Range("N73").Value = -6606.857400000001

# ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3154
$ws.Range("I76").Value = 2910.4443
$ws.Range("K76").Value = 2910.4443
$ws.Range("M76").Value = -2595.4443

# ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3154
$ws.Range("I79").Value = 2910.4443
$ws.Range("K79").Value = 2910.4443
$ws.Range("M79").Value = -1818.4443

# ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 736.6
$ws.Range("I98").Value = 736.6
$ws.Range("K98").Value = 736.6
$ws.Range("M98").Value = 761.4

# ALC row 108
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H108").Value = 99278
$ws.Range("J108").Value = 99278
$ws.Range("L108").Value = 99278
$ws.Range("N108").Value = -106958

# ALC row 110
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H110").Value = 55989
$ws.Range("J110").Value = 55989
$ws.Range("L110").Value = 55989
$ws.Range("N110").Value = -64169

# ALC row 111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 602.6667
$ws.Range("J111").Value = 912.6667
$ws.Range("L111").Value = 2738.0001
$ws.Range("N111").Value = -8872.000100000001

# ALC row 117
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H117").Value = 90738
$ws.Range("J117").Value = 90738
$ws.Range("L117").Value = 90738
$ws.Range("N117").Value = -99916

# ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 736.6
$ws.Range("I122").Value = 736.6
$ws.Range("K122").Value = 2209.8
$ws.Range("M122").Value = 240.1999999999998

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1377.8572
$ws.Range("I132").Value = 1155.8667
$ws.Range("K132").Value = 3467.6001
$ws.Range("M132").Value = -937.6001000000001

# ALC row 134
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H134").Value = 98978.664
$ws.Range("J134").Value = 98978.664
$ws.Range("L134").Value = 98978.664
$ws.Range("N134").Value = -109118.664

# ALC row 136
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H136").Value = 77977.336
$ws.Range("J136").Value = 77977.336
$ws.Range("L136").Value = 77977.336
$ws.Range("N136").Value = -88177.336

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2282.8
$ws.Range("I138").Value = 2104
$ws.Range("K138").Value = 6312
$ws.Range("M138").Value = -1172

# ARM row 104
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H104").Value = 33817.6
$ws.Range("J104").Value = 33817.6
$ws.Range("L104").Value = 33817.6
$ws.Range("N104").Value = -40805.6

# ARM row 107
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H107").Value = 52558.668
$ws.Range("J107").Value = 52558.668
$ws.Range("L107").Value = 52558.668
$ws.Range("N107").Value = -60238.668

# ARM row 121
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H121").Value = 86393.86
$ws.Range("J121").Value = 86393.86
$ws.Range("L121").Value = 86393.86
$ws.Range("N121").Value = -89887.86

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3036.5
$ws.Range("I122").Value = 2455.5715
$ws.Range("J122").Value = 3849.8
$ws.Range("K122").Value = 7366.7145
$ws.Range("L122").Value = 11549.4
$ws.Range("M122").Value = -4916.7145
$ws.Range("N122").Value = -16449.4

# BSM row 2
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 28998.666
$ws.Range("J2").Value = 28998.666
$ws.Range("L2").Value = 28998.666
$ws.Range("N2").Value = -29224.666

# BSM row 6
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2948.1667
$ws.Range("I20").Value = 2769.389
$ws.Range("K20").Value = 2769.389
$ws.Range("M20").Value = -2522.389

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1768.1154
$ws.Range("I94").Value = 1499.9412
$ws.Range("J94").Value = 2274.6667
$ws.Range("K94").Value = 1499.9412
$ws.Range("L94").Value = 2274.6667
$ws.Range("M94").Value = -1048.9412
$ws.Range("N94").Value = -3176.6667

# BSM row 110
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H110").Value = 37423.832
$ws.Range("J110").Value = 37423.832
$ws.Range("L110").Value = 37423.832
$ws.Range("N110").Value = -45603.832

# BSM row 119
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H119").Value = 98919.25
$ws.Range("J119").Value = 98919.25
$ws.Range("L119").Value = 98919.25
$ws.Range("N119").Value = -108595.25

# BSM row 132
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 29587.883
$ws.Range("J132").Value = 29587.883
$ws.Range("L132").Value = 29587.883
$ws.Range("N132").Value = -39707.883

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3310.9333
$ws.Range("I134").Value = 2133.9429
$ws.Range("J134").Value = 7430.4
$ws.Range("K134").Value = 6401.8287
$ws.Range("L134").Value = 22291.2
$ws.Range("M134").Value = -3866.8287
$ws.Range("N134").Value = -27361.2

# BSM row 135
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 109926.336
$ws.Range("J135").Value = 109926.336
$ws.Range("L135").Value = 109926.336
$ws.Range("N135").Value = -120066.336

# BSM row 138
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 77998.39999999999
$ws.Range("J138").Value = 77998.39999999999
$ws.Range("L138").Value = 77998.39999999999
$ws.Range("N138").Value = -88278.39999999999

# BSM row 140
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 43499.668
$ws.Range("J140").Value = 43499.668
$ws.Range("L140").Value = 43499.668
$ws.Range("N140").Value = -53859.668

# CRP row 9
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 79996
$ws.Range("J9").Value = 79996
$ws.Range("L9").Value = 79996
$ws.Range("N9").Value = -80332

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4008.762
$ws.Range("I31").Value = 2459.8572
$ws.Range("K31").Value = 2459.8572
$ws.Range("M31").Value = -2164.8572

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4008.762
$ws.Range("I34").Value = 2459.8572
$ws.Range("K34").Value = 2459.8572
$ws.Range("M34").Value = -2257.8572

# CRP row 108
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H108").Value = 29089.455
$ws.Range("J108").Value = 29089.455
$ws.Range("L108").Value = 29089.455
$ws.Range("N108").Value = -36769.455

# CRP row 114
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H114").Value = 62246
$ws.Range("J114").Value = 62246
$ws.Range("L114").Value = 62246
$ws.Range("N114").Value = -70924

# CRP row 138
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H138").Value = 53918.4
$ws.Range("J138").Value = 54898
$ws.Range("L138").Value = 54898
$ws.Range("N138").Value = -65178

# CUL row 114
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 5390.5293
$ws.Range("I114").Value = 614.5
$ws.Range("J114").Value = 9635.888999999999
$ws.Range("K114").Value = 1843.5
$ws.Range("L114").Value = 28907.667
$ws.Range("M114").Value = 1410.5
$ws.Range("N114").Value = -35415.667

# CUL row 121
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 1788.6111
$ws.Range("J121").Value = 3033
$ws.Range("L121").Value = 9099
$ws.Range("N121").Value = -11719

# GSM row 108
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H108").Value = 51420.57
$ws.Range("J108").Value = 51420.57
$ws.Range("L108").Value = 51420.57
$ws.Range("N108").Value = -59100.57

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 717984.5600000001
$ws.Range("I122").Value = 913144.0600000001
$ws.Range("K122").Value = 2739432.18
$ws.Range("M122").Value = -2736982.18

# GSM row 135
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 25000
$ws.Range("J135").Value = 25000
$ws.Range("L135").Value = 25000
$ws.Range("N135").Value = -35140

# GSM row 140
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 94332
$ws.Range("J140").Value = 94970.664
$ws.Range("L140").Value = 94970.664
$ws.Range("N140").Value = -105330.664

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1913.8667
$ws.Range("I7").Value = 975.6667
$ws.Range("K7").Value = 975.6667
$ws.Range("M7").Value = -863.6667

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6175934.5
$ws.Range("I40").Value = 2470.8
$ws.Range("K40").Value = 2470.8
$ws.Range("M40").Value = -2334.8

# LTW row 43
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 32602.8
$ws.Range("J43").Value = 32602.8
$ws.Range("L43").Value = 32602.8
$ws.Range("N43").Value = -32988.8

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 63160850
$ws.Range("I122").Value = 83335496
$ws.Range("K122").Value = 250006488
$ws.Range("M122").Value = -250004038

# LTW row 123
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H123").Value = 70767
$ws.Range("J123").Value = 74163.164
$ws.Range("L123").Value = 74163.164
$ws.Range("N123").Value = -83963.164

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 1913.8667
$ws.Range("I126").Value = 975.6667
$ws.Range("K126").Value = 2927.0001
$ws.Range("M126").Value = -457.0001000000002

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2431.8235
$ws.Range("I136").Value = 2431.8235
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 7295.470499999999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -4745.470499999999
$ws.Range("N136").ClearContents()

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1474.7693
$ws.Range("I122").Value = 712.55554
$ws.Range("K122").Value = 2137.66662
$ws.Range("M122").Value = 312.33338

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 5099.5884
$ws.Range("I126").Value = 4406.2856
$ws.Range("J126").Value = 8335
$ws.Range("K126").Value = 13218.8568
$ws.Range("L126").Value = 25005
$ws.Range("M126").Value = -10748.8568
$ws.Range("N126").Value = -29945

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2387.9333
$ws.Range("I132").Value = 1959
$ws.Range("J132").Value = 3388.7778
$ws.Range("K132").Value = 5877
$ws.Range("L132").Value = 10166.3334
$ws.Range("M132").Value = -3347
$ws.Range("N132").Value = -15226.3334
